$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Locations")

# Rename / redisplay the existing "Sector 2" keyhole rooms to the Keyhole naming scheme
$ws.Range("D12").Value = "TraverseKeyholeSecond"
$ws.Range("E12").Value = "Keyhole / Second District"

$ws.Range("D13").Value = "TraverseKeyholeThird"
$ws.Range("E13").Value = "Keyhole / Third District"

# Fill in the previously empty rows (13-17 -> worksheet rows 15-19)
$ws.Range("B15").Value = "0x3"
$ws.Range("C15").Value = "0x1"
$ws.Range("D15").Value = "TraverseThird"
$ws.Range("E15").Value = "Third District"

$ws.Range("B16").Value = "0xa"
$ws.Range("C16").Value = "0x1"
$ws.Range("D16").Value = "TraverseSecondCorrupted"
$ws.Range("E16").Value = "Second District"

$ws.Range("B17").Value = "0x5"
$ws.Range("C17").Value = "0x1"
$ws.Range("D17").Value = "TraverseKeyholeFirst"
$ws.Range("E17").Value = "Keyhole / First District"

$ws.Range("B18").Value = "0x8"
$ws.Range("C18").Value = "0x1"
$ws.Range("D18").Value = "TraverseKeyholeFall"
$ws.Range("E18").Value = "Keyhole / Terminus"

$ws.Range("B19").Value = "0x9"
$ws.Range("C19").Value = "0x1"
$ws.Range("D19").Value = "TraverseKeyholeTerminus"
$ws.Range("E19").Value = "Keyhole / Terminus"

# Update the selected cell on the Locations sheet
$null = $ws.Range("B20").Select()
